$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Add the new indicator row (row 82) - Trade Balance
$newRow = 82

$ws.Cells.Item($newRow, 1).Value = "GEM"
$ws.Cells.Item($newRow, 2).Value = "TRDBAL"
$ws.Cells.Item($newRow, 3).Value = "AUSTRALI"
$ws.Cells.Item($newRow, 4).Formula = "=B82&`",`"&C82"
$ws.Cells.Item($newRow, 5).Value = "Trade Balance"
$ws.Cells.Item($newRow, 6).Value = "AUS$"
$ws.Cells.Item($newRow, 7).Value = "Avg"
$ws.Cells.Item($newRow, 8).Value = "GDP & Domestic Demand"

# Update the selection / view to match the saved state
$ws.Range("F83").Select()
